# Updated cryptos list (prices and 1h volume changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.294.28'
$ws.Range("D3").Value = '1.577.71'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '''208.08'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("D8").Value = '''22.28'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '1.801.90'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("D13").Value = '1.575.69'
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").Value = '''3.79'
$ws.Range("D15").Value = '''0.520'
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").Value = '27.299.20'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").Value = '''62.58'
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = '''215.52'
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '0.0₃0689'
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("D23").Value = '''9.44'
$ws.Range("E23").Value = '  -3.56%  '
$ws.Range("D24").Value = '''2.01'
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("D25").Value = '''151.83'
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("E26").Value = '  -4.26%  '
$ws.Range("D27").Value = '''14.98'
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("B28").Value = 'BinanceUSD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '''0.104'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("E31").Value = '  -2.03%  '
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("D33").Value = '1.411.54'
$ws.Range("E33").Value = '  +2.13%  '
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").Value = '''2.28'
$ws.Range("E36").Value = '  -1.71%  '
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("D39").Value = '''0.822'
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("E40").Value = '  -2.78%  '
$ws.Range("E42").Value = '  +1.96%  '
$ws.Range("D43").Value = '''1.81'
$ws.Range("E43").Value = '  +3.07%  '
$ws.Range("D44").Value = '''5.34'
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("D45").Value = '''64.07'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").Value = '1.713.72'
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("D48").Value = '''86.36'
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").Value = '0.0₇0991'
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("D50").Value = '''0.0955'
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("E51").Value = '  -0.42%  '
